$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.15368166566049
$ws.Cells.Item(2, 3).Value = 6.463522822095095
$ws.Cells.Item(2, 4).Value = 4.464487150825882
$ws.Cells.Item(2, 6).Value = 21.22672403823336
$ws.Cells.Item(2, 7).Value = 23.47869156624686
$ws.Cells.Item(2, 8).Value = 12.99162733252838
$ws.Cells.Item(2, 11).Value = 8.095678867533417
$ws.Cells.Item(2, 13).Value = 19.97100899796924
$ws.Cells.Item(2, 14).Value = 17.79734019690946
$ws.Cells.Item(2, 15).Value = 19.06319147766586
$ws.Cells.Item(3, 2).Value = 6.888810264094829
$ws.Cells.Item(3, 3).Value = 6.439609489212041
$ws.Cells.Item(3, 4).Value = 4.37717067640287
$ws.Cells.Item(3, 6).Value = 21.23777867201598
$ws.Cells.Item(3, 7).Value = 23.49262747816026
$ws.Cells.Item(3, 8).Value = 13.0285557931943
$ws.Cells.Item(3, 11).Value = 7.912191942436116
$ws.Cells.Item(3, 13).Value = 19.37077490633287
$ws.Cells.Item(3, 14).Value = 17.85857551848901
$ws.Cells.Item(3, 15).Value = 19.11632312633718
$ws.Cells.Item(4, 2).Value = 6.721759395356923
$ws.Cells.Item(4, 3).Value = 6.424922149331342
$ws.Cells.Item(4, 4).Value = 4.321872831477076
$ws.Cells.Item(4, 6).Value = 21.2499662873048
$ws.Cells.Item(4, 7).Value = 23.50866549928375
$ws.Cells.Item(4, 8).Value = 13.05306866727468
$ws.Cells.Item(4, 11).Value = 7.795945484374031
$ws.Cells.Item(4, 13).Value = 19.00070946132391
$ws.Cells.Item(4, 14).Value = 17.89781554289431
$ws.Cells.Item(4, 15).Value = 19.15272519175
$ws.Cells.Item(5, 2).Value = 6.652679993185086
$ws.Cells.Item(5, 3).Value = 6.418938624150393
$ws.Cells.Item(5, 4).Value = 4.298930389214157
$ws.Cells.Item(5, 6).Value = 21.25628993332601
$ws.Cells.Item(5, 7).Value = 23.51707793819088
$ws.Cells.Item(5, 8).Value = 13.06352032926234
$ws.Cells.Item(5, 11).Value = 7.747711216393076
$ws.Cells.Item(5, 13).Value = 18.8497774233734
$ws.Cells.Item(5, 14).Value = 17.91422052383347
$ws.Cells.Item(5, 15).Value = 19.16850798973366
$ws.Cells.Item(6, 2).Value = 6.641152044090738
$ws.Cells.Item(6, 3).Value = 6.417945247730482
$ws.Cells.Item(6, 4).Value = 4.295096637642075
$ws.Cells.Item(6, 6).Value = 21.25742189771304
$ws.Cells.Item(6, 7).Value = 23.51858801769017
$ws.Cells.Item(6, 8).Value = 13.06528375523451
$ws.Cells.Item(6, 11).Value = 7.739651010007511
$ws.Cells.Item(6, 13).Value = 18.82471526535531
$ws.Cells.Item(6, 14).Value = 17.91696963218736
$ws.Cells.Item(6, 15).Value = 19.17118595949545
$ws.Cells.Item(7, 2).Value = 6.720831686318503
$ws.Cells.Item(7, 3).Value = 6.424841442936696
$ws.Cells.Item(7, 4).Value = 4.3215650529054
$ws.Cells.Item(7, 6).Value = 21.25004607700539
$ws.Cells.Item(7, 7).Value = 23.50877136009456
$ws.Cells.Item(7, 8).Value = 13.05320774938853
$ws.Cells.Item(7, 11).Value = 7.795298421335915
$ws.Cells.Item(7, 13).Value = 18.99867410922743
$ws.Cells.Item(7, 14).Value = 17.89803510639278
$ws.Cells.Item(7, 15).Value = 19.15293420521297
$ws.Cells.Item(8, 2).Value = 7.063331547099059
$ws.Cells.Item(8, 3).Value = 6.455280116008454
$ws.Cells.Item(8, 4).Value = 4.434739303015622
$ws.Cells.Item(8, 6).Value = 21.22941464137032
$ws.Cells.Item(8, 7).Value = 23.48194164854124
$ws.Cells.Item(8, 8).Value = 13.00397872963487
$ws.Cells.Item(8, 11).Value = 8.033178933494771
$ws.Cells.Item(8, 13).Value = 19.76451236302207
$ws.Cells.Item(8, 14).Value = 17.81811447817203
$ws.Cells.Item(8, 15).Value = 19.08072599569626
$ws.Cells.Item(9, 2).Value = 7.695784997040795
$ws.Cells.Item(9, 3).Value = 6.514812007800325
$ws.Cells.Item(9, 4).Value = 4.64266025214767
$ws.Cells.Item(9, 6).Value = 21.23180843365896
$ws.Cells.Item(9, 7).Value = 23.48882917090387
$ws.Cells.Item(9, 8).Value = 12.9220240665207
$ws.Cells.Item(9, 11).Value = 8.469653738446377
$ws.Cells.Item(9, 13).Value = 21.24378072099007
$ws.Cells.Item(9, 14).Value = 17.67433712461317
$ws.Cells.Item(9, 15).Value = 18.96917967587667
$ws.Cells.Item(10, 2).Value = 8.131753797379911
$ws.Cells.Item(10, 3).Value = 6.558289445493118
$ws.Cells.Item(10, 4).Value = 4.786079835568147
$ws.Cells.Item(10, 6).Value = 21.25965393718564
$ws.Cells.Item(10, 7).Value = 23.53026760833403
$ws.Cells.Item(10, 8).Value = 12.87069461178326
$ws.Cells.Item(10, 11).Value = 8.770121060221586
$ws.Cells.Item(10, 13).Value = 22.30394082352323
$ws.Cells.Item(10, 14).Value = 17.57648788830871
$ws.Cells.Item(10, 15).Value = 18.90563821706733
$ws.Cells.Item(11, 2).Value = 8.323031912628442
$ws.Cells.Item(11, 3).Value = 6.577976731937552
$ws.Cells.Item(11, 4).Value = 4.849145931960342
$ws.Cells.Item(11, 6).Value = 21.27796190959796
$ws.Cells.Item(11, 7).Value = 23.55701111048948
$ws.Cells.Item(11, 8).Value = 12.84927066808941
$ws.Cells.Item(11, 11).Value = 8.902066728993848
$ws.Cells.Item(11, 13).Value = 22.77805642583348
$ws.Cells.Item(11, 14).Value = 17.53364045523684
$ws.Cells.Item(11, 15).Value = 18.88074680310076
$ws.Cells.Item(12, 2).Value = 8.394393433019856
$ws.Cells.Item(12, 3).Value = 6.585415726605684
$ws.Cells.Item(12, 4).Value = 4.872702967067174
$ws.Cells.Item(12, 6).Value = 21.28570237110786
$ws.Cells.Item(12, 7).Value = 23.56827010075706
$ws.Cells.Item(12, 8).Value = 12.84143484372137
$ws.Cells.Item(12, 11).Value = 8.951323114491073
$ws.Cells.Item(12, 13).Value = 22.9562434125703
$ws.Cells.Item(12, 14).Value = 17.51765283712699
$ws.Cells.Item(12, 15).Value = 18.87189951241364
$ws.Cells.Item(13, 2).Value = 8.379072945656603
$ws.Cells.Item(13, 3).Value = 6.583814373700275
$ws.Cells.Item(13, 4).Value = 4.867644158998027
$ws.Cells.Item(13, 6).Value = 21.28399946962905
$ws.Cells.Item(13, 7).Value = 23.56579500499846
$ws.Cells.Item(13, 8).Value = 12.84311011110672
$ws.Cells.Item(13, 11).Value = 8.94074679309154
$ws.Cells.Item(13, 13).Value = 22.91793035089509
$ws.Cells.Item(13, 14).Value = 17.52108550719424
$ws.Cells.Item(13, 15).Value = 18.87377918347601
$ws.Cells.Item(14, 2).Value = 8.328924645000004
$ws.Cells.Item(14, 3).Value = 6.578589083974727
$ws.Cells.Item(14, 4).Value = 4.85109055591397
$ws.Cells.Item(14, 6).Value = 21.27858255906726
$ws.Cells.Item(14, 7).Value = 23.55791473488351
$ws.Cells.Item(14, 8).Value = 12.84862045952074
$ws.Cells.Item(14, 11).Value = 8.90613342993521
$ws.Cells.Item(14, 13).Value = 22.79274399254749
$ws.Cells.Item(14, 14).Value = 17.53232038811019
$ws.Cells.Item(14, 15).Value = 18.88000732891155
$ws.Cells.Item(15, 2).Value = 8.298066247719142
$ws.Cells.Item(15, 3).Value = 6.575386247626976
$ws.Cells.Item(15, 4).Value = 4.840908371737507
$ws.Cells.Item(15, 6).Value = 21.27536960336314
$ws.Cells.Item(15, 7).Value = 23.55323511636696
$ws.Cells.Item(15, 8).Value = 12.85203177479438
$ws.Cells.Item(15, 11).Value = 8.884838700518147
$ws.Cells.Item(15, 13).Value = 22.7158829177574
$ws.Cells.Item(15, 14).Value = 17.53923299472571
$ws.Cells.Item(15, 15).Value = 18.88389762981891
$ws.Cells.Item(16, 2).Value = 8.119106389584806
$ws.Cells.Item(16, 3).Value = 6.557000743107549
$ws.Cells.Item(16, 4).Value = 4.781913499455661
$ws.Cells.Item(16, 6).Value = 21.25857069516145
$ws.Cells.Item(16, 7).Value = 23.52867844525665
$ws.Cells.Item(16, 8).Value = 12.87213347380027
$ws.Cells.Item(16, 11).Value = 8.761400373275471
$ws.Cells.Item(16, 13).Value = 22.27277655812842
$ws.Cells.Item(16, 14).Value = 17.57932142935929
$ws.Cells.Item(16, 15).Value = 18.90734581803445
$ws.Cells.Item(17, 2).Value = 8.007471763998048
$ws.Cells.Item(17, 3).Value = 6.545696370465633
$ws.Cells.Item(17, 4).Value = 4.745156015363595
$ws.Cells.Item(17, 6).Value = 21.24970763594655
$ws.Cells.Item(17, 7).Value = 23.51563350037117
$ws.Cells.Item(17, 8).Value = 12.88495852477909
$ws.Cells.Item(17, 11).Value = 8.684441241340984
$ws.Cells.Item(17, 13).Value = 21.99872427048197
$ws.Cells.Item(17, 14).Value = 17.6043396089157
$ws.Cells.Item(17, 15).Value = 18.92275958169579
$ws.Cells.Item(18, 2).Value = 7.942602204492986
$ws.Cells.Item(18, 3).Value = 6.539186031420492
$ws.Cells.Item(18, 4).Value = 4.723809898810504
$ws.Cells.Item(18, 6).Value = 21.24514105198658
$ws.Cells.Item(18, 7).Value = 23.50887360076551
$ws.Cells.Item(18, 8).Value = 12.89251644777309
$ws.Cells.Item(18, 11).Value = 8.639731856220743
$ws.Cells.Item(18, 13).Value = 21.84033839096909
$ws.Cells.Item(18, 14).Value = 17.61888618730378
$ws.Cells.Item(18, 15).Value = 18.93200290631201
$ws.Cells.Item(19, 2).Value = 7.920526981444942
$ws.Cells.Item(19, 3).Value = 6.536980400363715
$ws.Cells.Item(19, 4).Value = 4.716547775264876
$ws.Cells.Item(19, 6).Value = 21.24368621534619
$ws.Cells.Item(19, 7).Value = 23.50671253421589
$ws.Cells.Item(19, 8).Value = 12.89510656809511
$ws.Cells.Item(19, 11).Value = 8.624518545448794
$ws.Cells.Item(19, 13).Value = 21.78658711626607
$ws.Cells.Item(19, 14).Value = 17.62383838523457
$ws.Cells.Item(19, 15).Value = 18.93519737041699
$ws.Cells.Item(20, 2).Value = 8.019424241870386
$ws.Cells.Item(20, 3).Value = 6.546900628163413
$ws.Cells.Item(20, 4).Value = 4.749090132719902
$ws.Cells.Item(20, 6).Value = 21.25059616810337
$ws.Cells.Item(20, 7).Value = 23.51694525624186
$ws.Cells.Item(20, 8).Value = 12.88357451315258
$ws.Cells.Item(20, 11).Value = 8.692679883645814
$ws.Cells.Item(20, 13).Value = 22.02797742129225
$ws.Cells.Item(20, 14).Value = 17.6016601661363
$ws.Cells.Item(20, 15).Value = 18.92107965589273
$ws.Cells.Item(21, 2).Value = 8.343683917875392
$ws.Cells.Item(21, 3).Value = 6.580124341292046
$ws.Cells.Item(21, 4).Value = 4.855961653111736
$ws.Cells.Item(21, 6).Value = 21.2801517503759
$ws.Cells.Item(21, 7).Value = 23.5601986771475
$ws.Cells.Item(21, 8).Value = 12.84699442015161
$ws.Cells.Item(21, 11).Value = 8.916319649532674
$ws.Cells.Item(21, 13).Value = 22.82955223452921
$ws.Cells.Item(21, 14).Value = 17.52901399089791
$ws.Cells.Item(21, 15).Value = 18.87816225984297
$ws.Cells.Item(22, 2).Value = 8.549339993574389
$ws.Cells.Item(22, 3).Value = 6.601742368192726
$ws.Cells.Item(22, 4).Value = 4.923910745903346
$ws.Cells.Item(22, 6).Value = 21.3041732965285
$ws.Cells.Item(22, 7).Value = 23.59506179088318
$ws.Cells.Item(22, 8).Value = 12.82470153147179
$ws.Cells.Item(22, 11).Value = 9.058340041923961
$ws.Cells.Item(22, 13).Value = 23.34548705790626
$ws.Cells.Item(22, 14).Value = 17.48292070309914
$ws.Cells.Item(22, 15).Value = 18.85348601813035
$ws.Cells.Item(23, 2).Value = 8.440168104581877
$ws.Cells.Item(23, 3).Value = 6.590214189251991
$ws.Cells.Item(23, 4).Value = 4.887822343432078
$ws.Cells.Item(23, 6).Value = 21.2909233918848
$ws.Cells.Item(23, 7).Value = 23.5758527322302
$ws.Cells.Item(23, 8).Value = 12.83645195987818
$ws.Cells.Item(23, 11).Value = 8.982928389563821
$ws.Cells.Item(23, 13).Value = 23.07090361840551
$ws.Cells.Item(23, 14).Value = 17.50739534420038
$ws.Cells.Item(23, 15).Value = 18.86634717043238
$ws.Cells.Item(24, 2).Value = 8.014022669232594
$ws.Cells.Item(24, 3).Value = 6.546356218235681
$ws.Cells.Item(24, 4).Value = 4.747312183140437
$ws.Cells.Item(24, 6).Value = 21.25019281503115
$ws.Cells.Item(24, 7).Value = 23.51634990683644
$ws.Cells.Item(24, 8).Value = 12.88419964972596
$ws.Cells.Item(24, 11).Value = 8.688956638356927
$ws.Cells.Item(24, 13).Value = 22.01475464208092
$ws.Cells.Item(24, 14).Value = 17.60287103347229
$ws.Cells.Item(24, 15).Value = 18.9218379613968
$ws.Cells.Item(25, 2).Value = 7.52941910521598
$ws.Cells.Item(25, 3).Value = 6.498740886555765
$ws.Cells.Item(25, 4).Value = 4.587998055298962
$ws.Cells.Item(25, 6).Value = 21.22657391717259
$ws.Cells.Item(25, 7).Value = 23.48057951740819
$ws.Cells.Item(25, 8).Value = 12.94263452596094
$ws.Cells.Item(25, 11).Value = 8.354995896025075
$ws.Cells.Item(25, 13).Value = 20.84741108593384
$ws.Cells.Item(25, 14).Value = 17.71185789988022
$ws.Cells.Item(25, 15).Value = 18.99612874715582
